# Edward Webb, Personal Project Logbook Journal (SDD)
# Add three new logbook entries after the "6.32pm 4/06/2017" entry, and
# make sure the trailing _GoBack bookmark ends up attached to the end of
# the new final paragraph (matching Word's normal behaviour of keeping
# _GoBack anchored at the last edit location).

$d = $word.ActiveDocument

$openQuote = [char]0x201C
$closeQuote = [char]0x201D

# Text of the three new logbook entries (paragraphs), each preceded by a
# paragraph break so that they become their own <w:p> elements.
$newText = "`r7.03pm 4/06/2017 – Added online help and added help button to program." + `
    "`r9.58pm 5/06/2017 – Received a 1000-word report from user for feedback." + `
    "`r11.37am 7/06/2017 – Set up Google Forms survey so users can give feedback. Sent it to various users, and put responses in ${openQuote}Survey Responses${closeQuote} folder in the project folder."

# Insert right at the end of the last paragraph's text (i.e. immediately
# before its paragraph mark / the _GoBack bookmark that currently lives
# there), so the existing "6.32pm ..." paragraph is left untouched and the
# new text becomes new paragraphs following it.
$lastPara = $d.Paragraphs.Last
$insertPos = $lastPara.Range.End
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertBefore($newText)

# After the insert above, the pre-existing "_GoBack" bookmark (originally
# sitting right before the old paragraph mark) ends up re-anchored in the
# middle of the document instead of at the end of the new last paragraph.
# Relocate it to sit right after the text of the new last paragraph, as in
# the target document.

# Append a temporary one-character marker at the very end of the document.
# (Adding a *collapsed* bookmark exactly at the document's absolute end
# position is unreliable, so we keep one extra character past the desired
# bookmark location while we re-create the bookmark, then remove it.)
$endPos = $d.Content.End
$marker = $d.Range($endPos, $endPos)
$marker.InsertBefore("X")

# Remove the old bookmark and re-add it immediately before the temporary
# marker character (i.e. exactly where it should end up once the marker is
# deleted).
$d.Bookmarks("_GoBack").Delete()
$bmPos = $d.Content.End - 2
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange)

# Remove the temporary marker character.
$delStart = $d.Content.End - 2
$delRange = $d.Range($delStart, $d.Content.End - 1)
$delRange.Delete()
